# "Two new Bottleneck pptx and excel"
# Adds the results for two newly-closed bottlenecks (f_numPP, f_numP) and
# removes the old unfilled "Function f_KK line line 534" placeholder row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9 - e_F4 / Function f_numPP line 673: fill in the solution/verification columns.
$ws.Range("C9").Value = "17th March, 2025"
$ws.Range("D9").Value = "23rd March, 2025"
$ws.Range("E9").Value = "Removed the unnecessary numPP variable."
$ws.Range("F9").Value = "Yes"
$ws.Range("H9").Value = "The function runtime has decreased by approximately 29%"
$ws.Range("H9").WrapText = $true

# Row 10 - Function f_numP line 636: fill in the solution/verification columns.
$ws.Range("C10").Value = "20th March, 2025"
$ws.Range("D10").Value = "25th March, 2025"
$ws.Range("E10").Value = "Removed the unnecessary numP variable."
$ws.Range("F10").Value = "Yes"
$ws.Range("H10").Value = "Function runtime has negligible runtime compared to program runtime(Almost 95% decrease)."
$ws.Range("H10").WrapText = $true
$ws.Rows(10).RowHeight = 43.5

# Remove the old, never-filled-in "Function f_KK line line 534" row entirely;
# everything below shifts up one row.
$ws.Rows(11).Delete()

# The sheet view no longer pins a frozen top-left cell; the new selection is
# the blank spacer row that used to be row 12.
$null = $ws.Range("A11:XFD11").Select()
